# Update computed profit figures across the Leve-crafting profit sheets
# (currentAveragePrice / Leve price / profit columns), per scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H17").Value = 962.625
$ws.Range("J17").Value = 1073.2307
$ws.Range("L17").Value = 3219.6921
$ws.Range("N17").Value = -3555.6921
$ws.Range("H100").Value = 2053.8333
$ws.Range("I100").Value = 1808.3334
$ws.Range("J100").Value = 2299.3333
$ws.Range("K100").Value = 1808.3334
$ws.Range("L100").Value = 2299.3333
$ws.Range("M100").Value = -1267.3334
$ws.Range("N100").Value = -3381.3333
$ws.Range("H107").Value = 611.125
$ws.Range("I107").Value = 614.75
$ws.Range("J107").Value = 600.25
$ws.Range("K107").Value = 614.75
$ws.Range("L107").Value = 600.25
$ws.Range("M107").Value = 1305.25
$ws.Range("N107").Value = -4440.25
$ws.Range("H110").Value = 35407.5
$ws.Range("J110").Value = 35407.5
$ws.Range("L110").Value = 35407.5
$ws.Range("N110").Value = -43587.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 3000300.8
$ws.Range("I11").Value = 3000300.8
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 3000300.8
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -3000156.8
$ws.Range("N11").ClearContents()
$ws.Range("H62").Value = 19500
$ws.Range("J62").Value = 19500
$ws.Range("L62").Value = 19500
$ws.Range("N62").Value = -20748
$ws.Range("H65").Value = 19500
$ws.Range("J65").Value = 19500
$ws.Range("L65").Value = 58500
$ws.Range("N65").Value = -64740
$ws.Range("H74").Value = 1379.7646
$ws.Range("I74").Value = 1217.7391
$ws.Range("J74").Value = 1718.5454
$ws.Range("K74").Value = 1217.7391
$ws.Range("L74").Value = 1718.5454
$ws.Range("M74").Value = -343.7391
$ws.Range("N74").Value = -3466.5454
$ws.Range("H77").Value = 1379.7646
$ws.Range("I77").Value = 1217.7391
$ws.Range("J77").Value = 1718.5454
$ws.Range("K77").Value = 6088.6955
$ws.Range("L77").Value = 8592.726999999999
$ws.Range("M77").Value = -1720.6955
$ws.Range("N77").Value = -17328.727
$ws.Range("H102").Value = 2358.6365
$ws.Range("I102").Value = 1843.6364
$ws.Range("J102").Value = 2873.6365
$ws.Range("K102").Value = 1843.6364
$ws.Range("L102").Value = 2873.6365
$ws.Range("M102").Value = -221.6364000000001
$ws.Range("N102").Value = -6117.636500000001
$ws.Range("H107").Value = 15431.143
$ws.Range("J107").Value = 15431.143
$ws.Range("L107").Value = 15431.143
$ws.Range("N107").Value = -23111.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 22440.75
$ws.Range("I26").Value = 10587.667
$ws.Range("K26").Value = 10587.667
$ws.Range("M26").Value = -10295.667
$ws.Range("H92").Value = 25666.666
$ws.Range("J92").Value = 25666.666
$ws.Range("L92").Value = 25666.666
$ws.Range("N92").Value = -30658.666
$ws.Range("H99").Value = 2009.1666
$ws.Range("I99").Value = 3035
$ws.Range("J99").Value = 1804
$ws.Range("K99").Value = 3035
$ws.Range("L99").Value = 1804
$ws.Range("M99").Value = -1537
$ws.Range("N99").Value = -4800
$ws.Range("H107").Value = 1681.6666
$ws.Range("I107").Value = 1644.5555
$ws.Range("J107").Value = 1904.3334
$ws.Range("K107").Value = 1644.5555
$ws.Range("L107").Value = 1904.3334
$ws.Range("M107").Value = 275.4445000000001
$ws.Range("N107").Value = -5744.3334
$ws.Range("H109").Value = 30945
$ws.Range("J109").Value = 30945
$ws.Range("L109").Value = 30945
$ws.Range("N109").Value = -33719

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1461.39
$ws.Range("I31").Value = 1044.909
$ws.Range("J31").Value = 2269.853
$ws.Range("K31").Value = 1044.909
$ws.Range("L31").Value = 2269.853
$ws.Range("M31").Value = -749.9090000000001
$ws.Range("N31").Value = -2859.853
$ws.Range("H34").Value = 1461.39
$ws.Range("I34").Value = 1044.909
$ws.Range("J34").Value = 2269.853
$ws.Range("K34").Value = 1044.909
$ws.Range("L34").Value = 2269.853
$ws.Range("M34").Value = -842.9090000000001
$ws.Range("N34").Value = -2673.853
$ws.Range("H43").Value = 19202.334
$ws.Range("J43").Value = 19202.334
$ws.Range("L43").Value = 19202.334
$ws.Range("N43").Value = -19570.334
$ws.Range("H92").Value = 30601
$ws.Range("J92").Value = 30601
$ws.Range("L92").Value = 30601
$ws.Range("N92").Value = -35593
$ws.Range("H101").Value = 19202.334
$ws.Range("J101").Value = 19202.334
$ws.Range("L101").Value = 19202.334
$ws.Range("N101").Value = -25692.334
$ws.Range("H105").Value = 1631.4
$ws.Range("I105").Value = 1618.5714
$ws.Range("J105").Value = 1811
$ws.Range("K105").Value = 1618.5714
$ws.Range("L105").Value = 1811
$ws.Range("M105").Value = 128.4286
$ws.Range("N105").Value = -5305
$ws.Range("H107").Value = 1572.7858
$ws.Range("I107").Value = 445.57144
$ws.Range("J107").Value = 2700
$ws.Range("K107").Value = 445.57144
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = 1474.42856
$ws.Range("N107").Value = -6540
$ws.Range("H134").Value = 1926.1666
$ws.Range("I134").Value = 1463.4
$ws.Range("J134").Value = 2256.7144
$ws.Range("K134").Value = 4390.200000000001
$ws.Range("L134").Value = 6770.1432
$ws.Range("M134").Value = -1855.200000000001
$ws.Range("N134").Value = -11840.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 718.14813
$ws.Range("I5").Value = 432.31818
$ws.Range("K5").Value = 1296.95454
$ws.Range("M5").Value = -1184.95454
$ws.Range("H40").Value = 4538
$ws.Range("I40").Value = 102.85714
$ws.Range("J40").Value = 20061
$ws.Range("K40").Value = 411.42856
$ws.Range("L40").Value = 80244
$ws.Range("M40").Value = -342.42856
$ws.Range("N40").Value = -80382
$ws.Range("H69").Value = 1006
$ws.Range("I69").Value = 741.3333
$ws.Range("J69").Value = 1800
$ws.Range("K69").Value = 2223.9999
$ws.Range("L69").Value = 5400
$ws.Range("N69").Value = -7022
$ws.Range("M69").Value = -1412.9999
$ws.Range("H72").Value = 1006
$ws.Range("I72").Value = 741.3333
$ws.Range("J72").Value = 1800
$ws.Range("K72").Value = 6671.9997
$ws.Range("L72").Value = 16200
$ws.Range("N72").Value = -24312
$ws.Range("M72").Value = -2615.9997
$ws.Range("H135").Value = 718.14813
$ws.Range("I135").Value = 432.31818
$ws.Range("K135").Value = 3890.86362
$ws.Range("M135").Value = -1355.86362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2844.2856
$ws.Range("I80").Value = 2726
$ws.Range("J80").Value = 3002
$ws.Range("K80").Value = 2726
$ws.Range("L80").Value = 3002
$ws.Range("M80").Value = -1728
$ws.Range("N80").Value = -4998
$ws.Range("H83").Value = 2844.2856
$ws.Range("I83").Value = 2726
$ws.Range("J83").Value = 3002
$ws.Range("K83").Value = 13630
$ws.Range("L83").Value = 15010
$ws.Range("M83").Value = -8638
$ws.Range("N83").Value = -24994
$ws.Range("H97").Value = 865.125
$ws.Range("I97").Value = 842
$ws.Range("J97").Value = 903.6667
$ws.Range("K97").Value = 842
$ws.Range("L97").Value = 903.6667
$ws.Range("M97").Value = -346
$ws.Range("N97").Value = -1895.6667
$ws.Range("H132").Value = 2565.7446
$ws.Range("I132").Value = 1486.84
$ws.Range("J132").Value = 3791.7727
$ws.Range("K132").Value = 4460.52
$ws.Range("L132").Value = 11375.3181
$ws.Range("M132").Value = -1930.52
$ws.Range("N132").Value = -16435.3181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060
$ws.Range("H136").Value = 2288.5715
$ws.Range("I136").Value = 1942.5
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 5827.5
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -3277.5
$ws.Range("N136").Value = -13350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 715
$ws.Range("I107").Value = 447.5
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 1342.5
$ws.Range("L107").Value = 3750
$ws.Range("M107").Value = 577.5
$ws.Range("N107").Value = -7590
